# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (rows 2-22) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 297
$ws1.Range("F3").Value  = 1189
$ws1.Range("F4").Value  = 16791
$ws1.Range("F6").Value  = 1644
$ws1.Range("F7").Value  = 66
$ws1.Range("F10").Value = 217
$ws1.Range("F12").Value = 11648
$ws1.Range("F14").Value = 1324
$ws1.Range("F15").Value = 4619
$ws1.Range("F16").Value = 444
$ws1.Range("F19").Value = 892
$ws1.Range("F20").Value = 337
$ws1.Range("F22").Value = 5213

# --- Sheet: 全部类型 (rows 2-25) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 297
$ws4.Range("F4").Value  = 1189
$ws4.Range("F5").Value  = 16791
$ws4.Range("F7").Value  = 1644
$ws4.Range("F8").Value  = 66
$ws4.Range("F11").Value = 217
$ws4.Range("F15").Value = 11648
$ws4.Range("F17").Value = 1324
$ws4.Range("F18").Value = 4619
$ws4.Range("F19").Value = 444
$ws4.Range("F22").Value = 892
$ws4.Range("F23").Value = 337
$ws4.Range("F25").Value = 5213
